# Atualiza gabarito VR (Verificacao Suplementar) e notas finais - calculo III
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3: aluno ficou "ausente" -> agora tirou 7 na VR
$ws.Range("G3").Value = 7

# G14: nota da VR corrigida de 2.7 para 8
$ws.Range("G14").Value = 8

# G18: aluno ficou "ausente" -> agora tirou 6 na VR
$ws.Range("G18").Value = 6

# G19: aluno ficou "ausente" -> agora tirou 4.6 na VR
$ws.Range("G19").Value = 4.6

# G29: nota da VR corrigida de 2.1 para 7
$ws.Range("G29").Value = 7

# G30: aluno ficou "ausente" -> agora tirou 4.6 na VR
$ws.Range("G30").Value = 4.6

# I32: aluno ficou "ausente" -> agora tirou 6 na APVA2
$ws.Range("I32").Value = 6

# Atualiza celula ativa/selecao para refletir o ultimo resultado editado
$ws.Range("N32").Select()
